# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held Strike# values; this regenerates
# that column with the correct K values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(4, 1, 2, 8, 4, 3, 3, 6, 8, 3, 1, 1, 8, 4, 2, 4, 5, 1, 3)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
